# Corrección información de diccionarios
# Fix two typos in the dictionary description texts.

$wb = $excel.ActiveWorkbook

# Sheet "rewiev-estado": correct "texto de la reseña o descrpcion del servicio"
# -> "texto de la reseña o descripcion del servicio" (cell C6)
$wsReview = $wb.Worksheets.Item("rewiev-estado")
$wsReview.Range("C6").Value = "texto de la reseña o descripcion del servicio"

# Sheet "metadata": correct "categoria del negcio" -> "categoria del negocio" (cell C8)
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Range("C8").Value = "categoria del negocio"
